$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.36 = 21371.05 pesos`n✅ 21371.05 pesos = 5.33 = 953.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update tasas numbers on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 186.5
$wsTasas.Range("O10").Value = 3985.7
$wsTasas.Range("N12").Value = 4011
$wsTasas.Range("O12").Value = 179
